$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("G4").Value = 2.18
$ws.Range("H4").Value = 2.92
$ws.Range("I4").Value = 3.6
$ws.Range("P4").Value = 1.5
$ws.Range("Q4").Value = 2.4
$ws.Range("U4").Value = 9.5
$ws.Range("V4").Value = 9
$ws.Range("W4").Value = 21
$ws.Range("X4").Value = 20
$ws.Range("AA4").Value = 5.7
$ws.Range("AD4").Value = 7.9
$ws.Range("AE4").Value = 17.5
$ws.Range("AF4").Value = 13
$ws.Range("AG4").Value = 55

# Row 6
$ws.Range("H6").Value = 3.5
$ws.Range("I6").Value = 5.7
$ws.Range("P6").Value = 1.5
$ws.Range("Q6").Value = 2.27
$ws.Range("W6").Value = 11
$ws.Range("AA6").Value = 7.1
$ws.Range("AD6").Value = 11.75

# Row 7
$ws.Range("G7").Value = 1.85
$ws.Range("I7").Value = 3.8
$ws.Range("R7").Value = 1.72
$ws.Range("S7").Value = 1.9
$ws.Range("T7").Value = 7.3
$ws.Range("U7").Value = 8.75
$ws.Range("W7").Value = 15.5
$ws.Range("X7").Value = 14.5
$ws.Range("Y7").Value = 26
$ws.Range("Z7").Value = 10.5
$ws.Range("AB7").Value = 14.5
$ws.Range("AC7").Value = 65
$ws.Range("AD7").Value = 11.5
$ws.Range("AE7").Value = 21
$ws.Range("AF7").Value = 13
$ws.Range("AG7").Value = 55
$ws.Range("AH7").Value = 35
$ws.Range("AI7").Value = 40
$ws.Range("AJ7").Value = 500

# Row 11
$ws.Range("H11").Value = 4.9
$ws.Range("I11").Value = 6.1
$ws.Range("L11").Value = 1.11
$ws.Range("N11").Value = 1.35
$ws.Range("O11").Value = 2.95
$ws.Range("T11").Value = 12.5
$ws.Range("U11").Value = 10
$ws.Range("W11").Value = 11.5
$ws.Range("AA11").Value = 11
$ws.Range("AB11").Value = 15
$ws.Range("AC11").Value = 45
$ws.Range("AD11").Value = 27
$ws.Range("AE11").Value = 45
$ws.Range("AF11").Value = 19.5
$ws.Range("AI11").Value = 40

# Row 13
$ws.Range("K13").Value = 8

# Row 14
$ws.Range("H14").Value = 3.2
$ws.Range("I14").Value = 4.6
$ws.Range("L14").Value = 1.52
$ws.Range("M14").Value = 2.22
$ws.Range("N14").Value = 2.45
$ws.Range("O14").Value = 1.42
$ws.Range("P14").Value = 1.57
$ws.Range("Q14").Value = 2.12
$ws.Range("R14").Value = 2.27
$ws.Range("T14").Value = 4.85
$ws.Range("U14").Value = 6.7
$ws.Range("X14").Value = 19
$ws.Range("Y14").Value = 50
$ws.Range("Z14").Value = 6.4
$ws.Range("AB14").Value = 24
$ws.Range("AF14").Value = 17.5
$ws.Range("AG14").Value = 90
$ws.Range("AH14").Value = 65
